$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 6 de Abril de 2020 a las 22:52"

# Update per-country COVID-19 stats (columns B-H) to the newer snapshot
$ws.Range("B4").Value = 362326
$ws.Range("C4").Value = 25653
$ws.Range("D4").Value = 19313
$ws.Range("E4").Value = 332299
$ws.Range("F4").Value = 8871
$ws.Range("G4").Value = 1098
$ws.Range("H4").Value = 10714

$ws.Range("B52").Value = 1579
$ws.Range("C52").Value = 94
$ws.Range("D52").Value = 88
$ws.Range("E52").Value = 1445
$ws.Range("F52").Value = 50
$ws.Range("G52").Value = 11
$ws.Range("H52").Value = 46

$ws.Range("B53").Value = 1562
$ws.Range("C53").Value = 76
$ws.Range("D53").Value = 460
$ws.Range("E53").Value = 1096
$ws.Range("F53").Value = 11
$ws.Range("G53").Value = 2
$ws.Range("H53").Value = 6

$ws.Range("B54").Value = 1554
$ws.Range("C54").Value = 0
$ws.Range("D54").Value = 325
$ws.Range("E54").Value = 1181
$ws.Range("F54").Value = 94
$ws.Range("G54").Value = 2
$ws.Range("H54").Value = 48

$ws.Range("B101").Value = 253
$ws.Range("C101").Value = 69
$ws.Range("D101").Value = 26
$ws.Range("E101").Value = 217
$ws.Range("F101").Value = 0
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 10

$ws.Range("B102").Value = 245
$ws.Range("C102").Value = 4
$ws.Range("D102").Value = 95
$ws.Range("E102").Value = 150
$ws.Range("F102").Value = 8
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 0

$ws.Range("B103").Value = 244
$ws.Range("C103").Value = 17
$ws.Range("D103").Value = 7
$ws.Range("E103").Value = 230
$ws.Range("F103").Value = 2
$ws.Range("G103").Value = 0
$ws.Range("H103").Value = 7

$ws.Range("B104").Value = 241
$ws.Range("C104").Value = 14
$ws.Range("D104").Value = 5
$ws.Range("E104").Value = 236
$ws.Range("F104").Value = 3
$ws.Range("G104").Value = 0
$ws.Range("H104").Value = 0

$ws.Range("B105").Value = 238
$ws.Range("C105").Value = 6
$ws.Range("D105").Value = 35
$ws.Range("E105").Value = 198
$ws.Range("F105").Value = 2
$ws.Range("G105").Value = 0
$ws.Range("H105").Value = 5

$ws.Range("B106").Value = 233
$ws.Range("C106").Value = 19
$ws.Range("D106").Value = 1
$ws.Range("E106").Value = 230
$ws.Range("F106").Value = 4
$ws.Range("G106").Value = 0
$ws.Range("H106").Value = 2

$ws.Range("B107").Value = 226
$ws.Range("C107").Value = 4
$ws.Range("D107").Value = 92
$ws.Range("E107").Value = 132
$ws.Range("F107").Value = 1
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = 2

$ws.Range("B108").Value = 216
$ws.Range("C108").Value = 69
$ws.Range("D108").Value = 33
$ws.Range("E108").Value = 179
$ws.Range("F108").Value = 5
$ws.Range("G108").Value = 3
$ws.Range("H108").Value = 4

$ws.Range("B109").Value = 214
$ws.Range("C109").Value = 0
$ws.Range("D109").Value = 31
$ws.Range("E109").Value = 178
$ws.Range("F109").Value = 2
$ws.Range("G109").Value = 0
$ws.Range("H109").Value = 5

$ws.Range("B110").Value = 188
$ws.Range("C110").Value = 14
$ws.Range("D110").Value = 39
$ws.Range("E110").Value = 147
$ws.Range("F110").Value = 6
$ws.Range("G110").Value = 0
$ws.Range("H110").Value = 2

$ws.Range("B117").Value = 151
$ws.Range("C117").Value = 2
$ws.Range("D117").Value = 50
$ws.Range("E117").Value = 97
$ws.Range("F117").Value = 20
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 4

$ws.Range("B172").Value = 14
$ws.Range("C172").Value = 0
$ws.Range("D172").Value = 0
$ws.Range("E172").Value = 14
$ws.Range("F172").Value = 0
$ws.Range("G172").Value = 0
$ws.Range("H172").Value = 0

$ws.Range("B173").Value = 14
$ws.Range("C173").Value = 2
$ws.Range("D173").Value = 0
$ws.Range("E173").Value = 14
$ws.Range("F173").Value = 0
$ws.Range("G173").Value = 0
$ws.Range("H173").Value = 0

$ws.Range("B177").Value = 12
$ws.Range("C177").Value = 1
$ws.Range("D177").Value = 0
$ws.Range("E177").Value = 12
$ws.Range("F177").Value = 0
$ws.Range("G177").Value = 0
$ws.Range("H177").Value = 0

$ws.Range("B178").Value = 12
$ws.Range("C178").Value = 0
$ws.Range("D178").Value = 0
$ws.Range("E178").Value = 12
$ws.Range("F178").Value = 2
$ws.Range("G178").Value = 0
$ws.Range("H178").Value = 0

$ws.Range("B184").Value = 10
$ws.Range("C184").Value = 1
$ws.Range("D184").Value = 0
$ws.Range("E184").Value = 9
$ws.Range("F184").Value = 0
$ws.Range("G184").Value = 0
$ws.Range("H184").Value = 1

$ws.Range("B185").Value = 10
$ws.Range("C185").Value = 0
$ws.Range("D185").Value = 0
$ws.Range("E185").Value = 9
$ws.Range("F185").Value = 0
$ws.Range("G185").Value = 0
$ws.Range("H185").Value = 1

$ws.Range("B186").Value = 10
$ws.Range("C186").Value = 1
$ws.Range("D186").Value = 4
$ws.Range("E186").Value = 6
$ws.Range("F186").Value = 0
$ws.Range("G186").Value = 0
$ws.Range("H186").Value = 0

$ws.Range("B187").Value = 9
$ws.Range("C187").Value = 0
$ws.Range("D187").Value = 0
$ws.Range("E187").Value = 9
$ws.Range("F187").Value = 0
$ws.Range("G187").Value = 0
$ws.Range("H187").Value = 0

$ws.Range("B188").Value = 9
$ws.Range("C188").Value = 0
$ws.Range("D188").Value = 1
$ws.Range("E188").Value = 8
$ws.Range("F188").Value = 0
$ws.Range("G188").Value = 0
$ws.Range("H188").Value = 0

$ws.Range("B197").Value = 6
$ws.Range("C197").Value = 0
$ws.Range("D197").Value = 1
$ws.Range("E197").Value = 5
$ws.Range("F197").Value = 0
$ws.Range("G197").Value = 0
$ws.Range("H197").Value = 0

$ws.Range("B199").Value = 6
$ws.Range("C199").Value = 0
$ws.Range("D199").Value = 0
$ws.Range("E199").Value = 5
$ws.Range("F199").Value = 0
$ws.Range("G199").Value = 0
$ws.Range("H199").Value = 1
